# Apply updated symbol list values (cryptos.xlsx refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextCell 2 4 "330.15"
Set-TextCell 2 5 "6.90%"
Set-TextCell 2 7 "12"

Set-TextCell 3 4 "40.22"
Set-TextCell 3 5 "8.11%"
Set-TextCell 3 7 "12"

Set-TextCell 4 4 "5.457"
Set-TextCell 4 5 "6.67%"
Set-TextCell 4 7 "12"

Set-TextCell 5 4 "0.08126"
Set-TextCell 5 5 "3.77%"
Set-TextCell 5 7 "12"

Set-TextCell 6 4 "4.529"
Set-TextCell 6 5 "3.05%"
Set-TextCell 6 7 "12"

Set-TextCell 7 4 "8.668"
Set-TextCell 7 5 "4.92%"
Set-TextCell 7 7 "12"

Set-TextCell 8 4 "1.917"
Set-TextCell 8 5 "1.77%"
Set-TextCell 8 7 "12"

Set-TextCell 9 5 "1.17%"
Set-TextCell 9 7 "12"

Set-TextCell 10 4 "0.9493"
Set-TextCell 10 5 "3.13%"
Set-TextCell 10 7 "12"

Set-TextCell 11 4 "0.1325"
Set-TextCell 11 5 "23.72%"
Set-TextCell 11 7 "12"

Set-TextCell 12 4 "0.2006"
Set-TextCell 12 5 "5.87%"
Set-TextCell 12 7 "12"

Set-TextCell 13 4 "0.09206"
Set-TextCell 13 5 "2.80%"
Set-TextCell 13 7 "12"

Set-TextCell 14 4 "0.03560"
Set-TextCell 14 5 "7.47%"
Set-TextCell 14 7 "12"

Set-TextCell 15 4 "0.09601"
Set-TextCell 15 5 "0.27%"
Set-TextCell 15 7 "12"

Set-TextCell 16 4 "0.001324"
Set-TextCell 16 5 "-3.73%"
Set-TextCell 16 7 "12"

Set-TextCell 17 4 "0.006152"
Set-TextCell 17 5 "7.85%"
Set-TextCell 17 7 "12"

Set-TextCell 18 4 "3.368"
Set-TextCell 18 5 "-0.54%"
Set-TextCell 18 7 "12"

Set-TextCell 19 4 "0.3515"
Set-TextCell 19 5 "1.89%"
Set-TextCell 19 7 "12"

Set-TextCell 20 4 "7.237"
Set-TextCell 20 5 "13.78%"
Set-TextCell 20 7 "12"

Set-TextCell 21 4 "0.1332"
Set-TextCell 21 5 "1.40%"
Set-TextCell 21 7 "12"

Set-TextCell 22 4 "0.2451"
Set-TextCell 22 5 "-0.77%"
Set-TextCell 22 7 "12"

Set-TextCell 23 5 "1.82%"
Set-TextCell 23 7 "12"

Set-TextCell 24 4 "0.001225"
Set-TextCell 24 5 "2.50%"
Set-TextCell 24 7 "12"

Set-TextCell 25 4 "0.004331"
Set-TextCell 25 5 "1.45%"
Set-TextCell 25 7 "12"

Set-TextCell 26 5 "-14.23%"
Set-TextCell 26 7 "12"

Set-TextCell 27 5 "37.70%"
Set-TextCell 27 7 "12"

Set-TextCell 28 7 "12"

Set-TextCell 29 7 "12"

Set-TextCell 30 7 "12"

Set-TextCell 31 7 "12"

Set-TextCell 32 7 "12"

Set-TextCell 33 7 "12"

Set-TextCell 34 7 "12"

Set-TextCell 35 7 "12"

Set-TextCell 36 7 "12"

Set-TextCell 37 7 "12"

Set-TextCell 38 7 "12"

Set-TextCell 39 4 "0.02531"
Set-TextCell 39 5 "16.27%"
Set-TextCell 39 7 "12"

Set-TextCell 40 4 "0.05246"
Set-TextCell 40 5 "4.34%"
Set-TextCell 40 7 "12"

Set-TextCell 41 4 "0.007741"
Set-TextCell 41 5 "2.22%"
Set-TextCell 41 7 "12"

Set-TextCell 42 4 "0.1428"
Set-TextCell 42 5 "5.61%"
Set-TextCell 42 7 "12"

Set-TextCell 43 4 "0.009266"
Set-TextCell 43 5 "7.28%"
Set-TextCell 43 7 "12"

Set-TextCell 44 4 "0.002160"
Set-TextCell 44 5 "4.42%"
Set-TextCell 44 7 "12"

Set-TextCell 45 4 "0.01084"
Set-TextCell 45 5 "35.61%"
Set-TextCell 45 7 "12"

Set-TextCell 46 4 "0.00006591"
Set-TextCell 46 5 "1.27%"
Set-TextCell 46 7 "12"

Set-TextCell 47 5 "0.09%"
Set-TextCell 47 7 "12"

Set-TextCell 48 2 "CoinbaseStockToken"
Set-TextCell 48 3 "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
Set-TextCell 48 4 "0.002400"
Set-TextCell 48 5 "139.52%"
Set-TextCell 48 7 "12"

Set-TextCell 49 2 "BOLO"
Set-TextCell 49 3 "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
Set-TextCell 49 4 "0.003345"
Set-TextCell 49 5 "1.53%"
Set-TextCell 49 7 "12"

Set-TextCell 50 4 "0.00002101"
Set-TextCell 50 5 "0.09%"
Set-TextCell 50 7 "12"

Set-TextCell 51 4 "0.0002001"
Set-TextCell 51 5 "0.09%"
Set-TextCell 51 7 "12"
